# Add a new time-log entry row (row 28) to the sheet:
#   Date: 12/08/2023 (serial 45268)
#   Name of Task: Internship
#   Description: Contributed technical work by aiding in resolving
#                inconsistencies flagged by the system for employee calls
#
# Commit message: "Task: Completed daily operations, 8 hours, 12/08"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row index is one past the previous last row (row 27 -> row 28)
$newRow = 28

# Column A: date value (12/08/2023), stored as the underlying date serial
# number (45268) and formatted to match the rest of the date column.
$ws.Cells.Item($newRow, 1).Value = 45268
$ws.Cells.Item($newRow, 1).NumberFormat = "d-mmm"

# Column B: task name
$ws.Cells.Item($newRow, 2).Value = "Internship"

# Column C: description text (matches the existing shared string used by
# the preceding rows)
$ws.Cells.Item($newRow, 3).Value = "Contributed technical work by aiding in resolving inconsistencies flagged by the system for employee calls"

# Update the active selection to the next empty row in column C, mirroring
# the author's cursor position after typing the new entry.
$ws.Range("C29").Select()
